$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.734.68'
$ws.Range("E2").Value = '  +4.92%  '
$ws.Range("D3").Value = '3.091.61'
$ws.Range("E3").Value = '  +2.86%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '558.59'
$ws.Range("E5").Value = '  +3.02%  '
$ws.Range("D6").Value = '143.84'
$ws.Range("E6").Value = '  +9.35%  '
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '3.082.42'
$ws.Range("E8").Value = '  +2.80%  '
$ws.Range("D9").Value = '0.498'
$ws.Range("E9").Value = '  +1.83%  '
$ws.Range("D10").Value = '7.12'
$ws.Range("E10").Value = '  +17.08%  '
$ws.Range("E11").Value = '  +4.70%  '
$ws.Range("E12").Value = '  +4.05%  '
$ws.Range("E13").Value = '  +4.33%  '
$ws.Range("D14").Value = '35.28'
$ws.Range("E14").Value = '  +2.73%  '
$ws.Range("D15").Value = '3.600.71'
$ws.Range("E15").Value = '  +3.12%  '
$ws.Range("D16").Value = '64.804.81'
$ws.Range("E16").Value = '  +5.03%  '
$ws.Range("D17").Value = '3.095.24'
$ws.Range("E17").Value = '  +3.08%  '
$ws.Range("E18").Value = '  -0.55%  '
$ws.Range("D19").Value = '6.79'
$ws.Range("E19").Value = '  +2.77%  '
$ws.Range("D20").Value = '482.06'
$ws.Range("E20").Value = '  -0.36%  '
$ws.Range("D21").Value = '13.82'
$ws.Range("E21").Value = '  +4.59%  '
$ws.Range("D22").Value = '0.675'
$ws.Range("E22").Value = '  +1.25%  '
$ws.Range("D23").Value = '7.54'
$ws.Range("E23").Value = '  +8.77%  '
$ws.Range("D24").Value = '13.35'
$ws.Range("E24").Value = '  +11.97%  '
$ws.Range("D25").Value = "'80.80"
$ws.Range("E25").Value = '  -1.64%  '
$ws.Range("E26").Value = '  +0.01%  '
$ws.Range("E27").Value = '  +3.71%  '
$ws.Range("D28").Value = '8.18'
$ws.Range("E28").Value = '  +6.77%  '
$ws.Range("E29").Value = '  +7.94%  '
$ws.Range("E30").Value = '  +0.18%  '
$ws.Range("D31").Value = '26.03'
$ws.Range("E31").Value = '  +1.43%  '
$ws.Range("D32").Value = '1.15'
$ws.Range("E32").Value = '  +3.61%  '
$ws.Range("E33").Value = '  +6.00%  '
$ws.Range("D34").Value = '5.68'
$ws.Range("E34").Value = '  +1.47%  '
$ws.Range("E35").Value = '  +6.54%  '
$ws.Range("D36").Value = '54.79'
$ws.Range("E36").Value = '  +0.04%  '
$ws.Range("D37").Value = '466.87'
$ws.Range("E37").Value = '  +6.97%  '
$ws.Range("D38").Value = '0.0407'
$ws.Range("E38").Value = '  +6.93%  '
$ws.Range("D39").Value = '0.0823'
$ws.Range("E39").Value = '  +3.83%  '
$ws.Range("D40").Value = '2.89'
$ws.Range("E40").Value = '  +19.91%  '
$ws.Range("D41").Value = '3.005.01'
$ws.Range("E41").Value = '  -4.10%  '
$ws.Range("E42").Value = '  +2.52%  '
$ws.Range("E43").Value = '  -1.27%  '
$ws.Range("D44").Value = '28.32'
$ws.Range("E44").Value = '  +7.81%  '
$ws.Range("D45").Value = '0.258'
$ws.Range("E45").Value = '  +7.39%  '
$ws.Range("E47").Value = '  +8.51%  '
$ws.Range("E48").Value = '  +4.16%  '
$ws.Range("D49").Value = '118.31'
$ws.Range("E49").Value = '  +2.58%  '
$ws.Range("E50").Value = '  +6.73%  '
$ws.Range("E51").Value = '  +2.74%  '
